$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: could not find '$old'"
    }
}

Replace-Text "2024-09-19 Thursday" "2024-09-20 Friday"
Replace-Text "887÷9=98, 5" "467÷6=77, 5"
Replace-Text "755÷4=188, 3" "617÷9=68, 5"
Replace-Text "656÷2=328, 0" "124÷4=31, 0"
Replace-Text "193÷6=32, 1" "998÷6=166, 2"
Replace-Text "363÷4=90, 3" "184÷3=61, 1"
Replace-Text "215÷5=43, 0" "509÷8=63, 5"
Replace-Text "352÷9=39, 1" "289÷8=36, 1"
Replace-Text "165÷7=23, 4" "216÷9=24, 0"
Replace-Text "252÷8=31, 4" "982÷4=245, 2"
Replace-Text "615÷3=205, 0" "120÷8=15, 0"
Replace-Text "251÷4=62, 3" "829÷4=207, 1"
Replace-Text "319÷2=159, 1" "472÷9=52, 4"
Replace-Text "700÷2=350, 0" "689÷5=137, 4"
Replace-Text "428÷4=107, 0" "420÷8=52, 4"
Replace-Text "580÷7=82, 6" "669÷6=111, 3"
Replace-Text "777÷8=97, 1" "811÷8=101, 3"
Replace-Text "803÷7=114, 5" "216÷4=54, 0"
Replace-Text "603÷3=201, 0" "125÷9=13, 8"
Replace-Text "220÷6=36, 4" "758÷5=151, 3"
Replace-Text "711÷6=118, 3" "644÷7=92, 0"
Replace-Text "787÷9=87, 4" "829÷7=118, 3"
Replace-Text "261÷3=87, 0" "700÷7=100, 0"
Replace-Text "678÷8=84, 6" "435÷3=145, 0"
Replace-Text "535÷5=107, 0" "976÷4=244, 0"
Replace-Text "396÷7=56, 4" "342÷8=42, 6"
